$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the "first eval" SD-filtered QoQ forecast error table
# (ifoCAST full series evaluation)

$data = @{
    2  = @{ B = 0.3895586535830836; C = 0.3895586535830836; D = 0.1926588447844201; E = 0.4389292024739526; F = 0.2145128497976859; G = 9 }
    3  = @{ B = 0.2985262890045486; C = 0.3172922402961569; D = 0.1252830370555511; E = 0.3539534391068282; F = 0.1994532552043115; G = 11 }
    4  = @{ B = 0.2787209351834913; C = 0.3278067591540666; D = 0.1386282432761488; E = 0.3723281392483635; F = 0.2589153759884639; G = 11 }
    5  = @{ B = 0.3034897228953755; C = 0.3268934070108915; D = 0.1294052293206626; E = 0.35972938345465;   F = 0.2064647113882407; G = 8 }
    6  = @{ B = 0.3209127895838344; C = 0.3481513866787824; D = 0.1477021037193441; E = 0.3843203139561375; F = 0.2242915978162668; G = 9 }
    7  = @{ B = 0.2986844067393156; C = 0.3331691208960183; D = 0.1354274215115737; E = 0.3680046487635362; F = 0.2298199212609901; G = 8 }
    8  = @{ B = 0.3065635236640661; C = 0.335248462641951;  D = 0.1402267591879591; E = 0.3744686357867092; F = 0.2322781938489615; G = 7 }
    9  = @{ B = 0.323074747095869;  C = 0.3439312876808158; D = 0.1492574881197667; E = 0.3863385666999435; F = 0.2288235751843726; G = 7 }
    10 = @{ B = 0.3092991727900825; C = 0.3325527292829896; D = 0.1414838076039549; E = 0.3761433338555329; F = 0.2393162899682282; G = 5 }
    11 = @{ B = 0.2921799865119131; C = 0.3130099903630257; D = 0.1357437534307621; E = 0.3684341914518278; F = 0.2591643208793996; G = 4 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
